# Update cfb_weather.xlsx with Timestamp 2024-12-30T16:21:31.938712
#
# This refreshes the scraped weather/odds numbers on the "FBS" sheet for the
# most-recently-updated games, including a full swap of the "Iowa @ Missouri"
# and "Louisville @ Washington" rows (their data had been in the wrong rows),
# a few wind-direction corrections, and bumps the Timestamp column (AK) to
# the new run time on every data row. It also fixes a wind-direction value on
# the "Other" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FBS")

# --- Row 2: South Carolina @ Illinois ---
$ws.Range("O2").Value = 37.94
$ws.Range("P2").Value = 19.1
$ws.Range("R2").Value = 0.8
$ws.Range("U2").Value = 7.4

# --- Row 3: Alabama @ Michigan ---
$ws.Range("O3").Value = 34.16
$ws.Range("P3").Value = 13.5
$ws.Range("R3").Value = 3.2
$ws.Range("U3").Value = 1.4
$ws.Range("AB3").Value = 14
$ws.Range("AF3").Value = -0.5

# --- Row 4: now Louisville @ Washington (was Iowa @ Missouri) ---
$ws.Range("A4").Value = "Louisville @ Washington"
$ws.Range("B4").Value = "TUE 12/31"
$ws.Range("C4").Value = "11:00 AM"
$ws.Range("D4").Value = "High"
$ws.Range("E4").Value = "E-W"
$ws.Range("F4").Value = "Low"
$ws.Range("G4").Value = "x w"
$ws.Range("H4").Value = -126.90537163
$ws.Range("I4").Value = 51.96
$ws.Range("J4").Value = 57.6
$ws.Range("K4").Value = 5.1
$ws.Range("L4").Value = 1920
$ws.Range("N4").Value = "NW"
$ws.Range("O4").Value = 43.1
$ws.Range("P4").Value = 2.1
$ws.Range("Q4").Value = "NW"
$ws.Range("R4").Value = 0
$ws.Range("U4").Value = -3
$ws.Range("V4").Value = "47.6503235, -122.3015746"
$ws.Range("W4").Value = 49.5
$ws.Range("X4").Value = -110
$ws.Range("Y4").Value = 49.5
$ws.Range("Z4").Value = -115
$ws.Range("AA4").Value = 2.5
$ws.Range("AB4").Value = 1
$ws.Range("AF4").Value = 1.5

# --- Row 5: now Iowa @ Missouri (was Louisville @ Washington) ---
$ws.Range("A5").Value = "Iowa @ Missouri"
$ws.Range("B5").Value = "MON 12/30"
$ws.Range("C5").Value = "01:30 PM"
$ws.Range("D5").Value = "Mid"
$ws.Range("E5").Value = "NE-SW"
$ws.Range("F5").Value = "High"
$ws.Range("G5").Value = "E/W"
$ws.Range("H5").Value = -7.949203499999982
$ws.Range("I5").Value = 56.44
$ws.Range("J5").Value = 51.51
$ws.Range("K5").Value = 6.9
$ws.Range("L5").Value = 1927
$ws.Range("M5").Value = "WNW"
$ws.Range("O5").Value = 49.28
$ws.Range("P5").Value = 10.4
$ws.Range("Q5").Value = "WNW"
$ws.Range("U5").Value = 3.5
$ws.Range("V5").Value = "38.9358491, -92.3332009"
$ws.Range("W5").Value = 40.5
$ws.Range("X5").Value = -114
$ws.Range("Y5").Value = 56.5
$ws.Range("Z5").Value = -110
$ws.Range("AA5").Value = -2.5
$ws.Range("AB5").Value = -1
$ws.Range("AE5").Value = 0.3950617283950617
$ws.Range("AF5").Value = -1.5

# --- Row 13 ---
$ws.Range("AB13").Value = -13
$ws.Range("AF13").Value = -1

# --- Row 14 ---
$ws.Range("AB14").Value = 9
$ws.Range("AF14").Value = -1.5

# --- Wind-direction (wind_dir_fg) corrections on a handful of rows ---
$ws.Range("Q28").Value = "NW"
$ws.Range("Q29").Value = "NNW"
$ws.Range("Q30").Value = "NW"
$ws.Range("Q31").Value = "N"
$ws.Range("Q34").Value = "NW"

# --- Bump the Timestamp column (AK / column 37) on every data row ---
for ($r = 2; $r -le 34; $r++) {
    $ws.Cells.Item($r, 37).Value = "2024-12-30T16:21:31.938712"
}

# --- "Other" sheet: wind-direction correction ---
$wsOther = $wb.Worksheets.Item("Other")
$wsOther.Range("S4").Value = "SSE"
